$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text with new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.81 = 10381.49 pesos`n✅ 10381.49 pesos = 2.79 = 932.52 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" sheet rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 356.5
$wsTasas.Range("O10").Value = 3701
$wsTasas.Range("N12").Value = 3725
$wsTasas.Range("O12").Value = 334.6
